$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.443.81'
$ws.Range("E2").Value = '  +0.97%  '
$ws.Range("D3").Value = '3.624.26'
$ws.Range("E3").Value = '  +2.79%  '
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = "'603.16"
$ws.Range("E5").Value = '  -0.76%  '
$ws.Range("D6").Value = "'196.90"
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("E7").Value = '  -0.90%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("E9").Value = '  +6.89%  '
$ws.Range("D10").Value = "'0.644"
$ws.Range("E10").Value = '  -1.03%  '
$ws.Range("D11").Value = "'53.34"
$ws.Range("E11").Value = '  -0.97%  '
$ws.Range("E12").Value = '  +0.67%  '
$ws.Range("D13").Value = "'9.53"
$ws.Range("E13").Value = '  +0.04%  '
$ws.Range("D14").Value = '4.195.48'
$ws.Range("E14").Value = '  +2.78%  '
$ws.Range("D15").Value = "'607.13"
$ws.Range("E15").Value = '  +1.72%  '
$ws.Range("D16").Value = "'12.98"
$ws.Range("E16").Value = '  +1.03%  '
$ws.Range("D17").Value = '70.497.76'
$ws.Range("E17").Value = '  +0.80%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.624.20'
$ws.Range("E18").Value = '  +2.83%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").Value = "'18.99"
$ws.Range("E19").Value = '  -0.83%  '
$ws.Range("E20").Value = '  +1.08%  '
$ws.Range("E21").Value = '  +0.28%  '
$ws.Range("D22").Value = "'18.14"
$ws.Range("E22").Value = '  -0.58%  '
$ws.Range("D23").Value = "'5.29"
$ws.Range("E23").Value = '  -0.21%  '
$ws.Range("D24").Value = "'103.57"
$ws.Range("E24").Value = '  +1.25%  '
$ws.Range("E25").Value = '  -1.79%  '
$ws.Range("D26").Value = "'2.99"
$ws.Range("E26").Value = '  -6.28%  '
$ws.Range("D27").Value = "'10.60"
$ws.Range("E27").Value = '  -2.61%  '
$ws.Range("E28").Value = '  +0.81%  '
$ws.Range("D29").Value = "'33.74"
$ws.Range("E29").Value = '  +0.66%  '
$ws.Range("D30").Value = "'4.70"
$ws.Range("E30").Value = '  +9.35%  '
$ws.Range("D31").Value = "'7.19"
$ws.Range("E31").Value = '  +1.52%  '
$ws.Range("D32").Value = "'12.27"
$ws.Range("E32").Value = '  -1.62%  '
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("D34").Value = "'63.31"
$ws.Range("D35").Value = '0.0₃0884'
$ws.Range("E35").Value = '  +3.05%  '
$ws.Range("D36").Value = '3.956.19'
$ws.Range("E36").Value = '  +6.02%  '
$ws.Range("E37").Value = '  +0.20%  '
$ws.Range("D38").Value = "'3.07"
$ws.Range("E38").Value = '  -0.21%  '
$ws.Range("D39").Value = "'516.37"
$ws.Range("E39").Value = '  +5.79%  '
$ws.Range("D40").Value = "'0.389"
$ws.Range("E40").Value = '  -1.11%  '
$ws.Range("D41").Value = "'36.59"
$ws.Range("E41").Value = '  -0.11%  '
$ws.Range("E42").Value = '  -3.03%  '
$ws.Range("E43").Value = '  +2.23%  '
$ws.Range("D44").Value = "'0.0461"
$ws.Range("E44").Value = '  +1.04%  '
$ws.Range("D45").Value = "'3.50"
$ws.Range("E45").Value = '  +6.53%  '
$ws.Range("D46").Value = "'2.91"
$ws.Range("E46").Value = '  +3.17%  '
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("D48").Value = "'8.56"
$ws.Range("E49").Value = '  -0.33%  '
$ws.Range("E50").Value = '  +0.68%  '
$ws.Range("E51").Value = '  -0.05%  '
